$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

# Update Mean (H) and Std (I) columns for the rows affected by the preferences-order evaluation update
$updates = @(
    @{ Row = 26; H = 0.46866; I = 0.02628 }
    @{ Row = 27; H = 0.04901; I = 0.02163 }
    @{ Row = 28; H = 0.43274; I = 0.02437 }
    @{ Row = 29; H = 0.00036; I = 0.00177 }
    @{ Row = 30; H = 0.46828; I = 0.02473 }
    @{ Row = 31; H = 0.04973; I = 0.02281 }
    @{ Row = 32; H = 0.43102; I = 0.02617 }
    @{ Row = 33; H = 0.00108; I = 0.00293 }
    @{ Row = 34; H = 0.52584; I = 0.01579 }
    @{ Row = 35; H = 0.0018; I = 0.00441 }
    @{ Row = 36; H = 0.52646; I = 0.01573 }
    @{ Row = 37; H = 0.0018; I = 0.00441 }
    @{ Row = 38; H = 0.5254799999999999; I = 0.01597 }
    @{ Row = 39; H = 0.0018; I = 0.00441 }
    @{ Row = 40; H = 0.52637; I = 0.01567 }
    @{ Row = 41; H = 0.0018; I = 0.00441 }
    @{ Row = 66; H = 0.42873; I = 0.01877 }
    @{ Row = 67; H = 0.03027; I = 0.01818 }
    @{ Row = 68; H = 0.40932; I = 0.01968 }
    @{ Row = 69; H = 0.00108; I = 0.00293 }
    @{ Row = 70; H = 0.43056; I = 0.01735 }
    @{ Row = 71; H = 0.03207; I = 0.01803 }
    @{ Row = 72; H = 0.40721; I = 0.02083 }
    @{ Row = 73; H = 0.00216; I = 0.00461 }
    @{ Row = 74; H = 0.5001100000000001; I = 0.01509 }
    @{ Row = 75; H = 0.00216; I = 0.00461 }
    @{ Row = 76; H = 0.50614; I = 0.01429 }
    @{ Row = 77; H = 0.0036; I = 0.0051 }
    @{ Row = 78; H = 0.50004; I = 0.01512 }
    @{ Row = 79; H = 0.00216; I = 0.00461 }
    @{ Row = 80; H = 0.50597; I = 0.0143 }
    @{ Row = 81; H = 0.0036; I = 0.0051 }
    @{ Row = 106; H = 0.46129; I = 0.01767 }
    @{ Row = 107; H = 0.04685; I = 0.02101 }
    @{ Row = 108; H = 0.42804; I = 0.0173 }
    @{ Row = 109; H = 0.00216; I = 0.00461 }
    @{ Row = 110; H = 0.46068; I = 0.01907 }
    @{ Row = 111; H = 0.04541; I = 0.02139 }
    @{ Row = 112; H = 0.4258; I = 0.01662 }
    @{ Row = 113; H = 0.00324; I = 0.00563 }
    @{ Row = 114; H = 0.52014; I = 0.01407 }
    @{ Row = 115; H = 0.00468; I = 0.00726 }
    @{ Row = 116; H = 0.52173; I = 0.01344 }
    @{ Row = 117; H = 0.00541; I = 0.00764 }
    @{ Row = 118; H = 0.5204299999999999; I = 0.01434 }
    @{ Row = 119; H = 0.00468; I = 0.00726 }
    @{ Row = 120; H = 0.52137; I = 0.01323 }
    @{ Row = 121; H = 0.00505; I = 0.00724 }
    @{ Row = 146; H = 0.4259; I = 0.02026 }
    @{ Row = 147; H = 0.02306; I = 0.01531 }
    @{ Row = 148; H = 0.40089; I = 0.01984 }
    @{ Row = 149; H = 0.00144; I = 0.0033 }
    @{ Row = 150; H = 0.42748; I = 0.01997 }
    @{ Row = 151; H = 0.02306; I = 0.01488 }
    @{ Row = 152; H = 0.40123; I = 0.01981 }
    @{ Row = 153; H = 0.00396; I = 0.00678 }
    @{ Row = 154; H = 0.49026; I = 0.01312 }
    @{ Row = 155; H = 0.00144; I = 0.0033 }
    @{ Row = 156; H = 0.49629; I = 0.01255 }
    @{ Row = 157; H = 0.00396; I = 0.00678 }
    @{ Row = 158; H = 0.49074; I = 0.01334 }
    @{ Row = 159; H = 0.00144; I = 0.0033 }
    @{ Row = 160; H = 0.49677; I = 0.01301 }
    @{ Row = 161; H = 0.00432; I = 0.0068 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 8).Value = $u.H   # Column H = Mean
    $ws.Cells.Item($u.Row, 9).Value = $u.I   # Column I = Std
}
